# The workbook's vocabulary sheet was regenerated from the source Google
# Sheet. Net effect versus the previous export:
#   1. F45's note about which terms use datacite:Other was trimmed down to
#      just the ResourceTypeGeneral reference.
#   2. A new "datacite:OtherDescriptionType" term was inserted right before
#      the existing "datacite:ContributorType" section header (old row 118),
#      pushing everything from there down by one row.
#   3. A new "datacite:OtherContributorType" term was inserted right before
#      the existing "datacite:ResourceCreatorType" section header (old row
#      139, now row 140 after the first insert), pushing everything below
#      it down by one more row.
#   4. Two more of the sheet's trailing blank "datacite:" filler rows were
#      appended at the very end, extending the used range from A1:S160 to
#      A1:S162.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Trim the F45 note.
$ws.Range("F45").Value = "datacite:ResourceTypeGeneral"

# 2) Insert "datacite:OtherDescriptionType" before the old row 118
#    ("datacite:ContributorType" header), shifting rows 118+ down by one.
$ws.Rows(118).Insert()
$ws.Range("A118").Value = "datacite:OtherDescriptionType"
$ws.Range("B118").Value = "Other"
$ws.Range("F118").Value = "datacite:DateType"

# 3) Insert "datacite:OtherContributorType" before the old row 139
#    ("datacite:ResourceCreatorType" header), which is now at row 140
#    after the previous insert, shifting rows 140+ down by one.
$ws.Rows(140).Insert()
$ws.Range("A140").Value = "datacite:OtherContributorType"
$ws.Range("B140").Value = "Other"
$ws.Range("F140").Value = "datacite:ContributorType"

# 4) Append two more blank "datacite:" filler rows at the new end of the
#    sheet (rows 161 and 162), matching the existing filler rows above them.
$ws.Range("A161").Value = "datacite:"
$ws.Range("A162").Value = "datacite:"
